$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 243, shifting rows 243:355 down to 244:356.
$ws.Rows(243).Insert()

# Populate the newly inserted row 243 with its data.
$ws.Range("A243").Value = 9
$ws.Range("B243").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C243").Value = "Metropolitana"
$ws.Range("D243").Value = 44917
$ws.Range("E243").Value = 13
$ws.Range("F243").Value = 300000001
$ws.Range("G243").Value = "Rabanito"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 7000
$ws.Range("K243").Value = 3000
$ws.Range("L243").Value = 3000
$ws.Range("M243").Value = 3000
$ws.Range("N243").Value = '$/cien unidades (volumen en unidades)'
$ws.Range("O243").Value = "Provincia de Chacabuco"
$ws.Range("P243").Value = 30
$ws.Range("Q243").Value = 100
$ws.Range("R243").Value = "Hortaliza"
